# Daily attendance processing - 2025-11-28 05:26:44
#
# The "Recorded By" column (G) stores a comma-separated list of the
# users/processes that recorded a session (e.g. "dnasr281@gmail.com, System").
# For every row where "System" appears as one of the comma-separated
# entries (and there is more than one entry), reverse the order of the
# entries in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    if ($text -notlike "*,*") {
        continue
    }

    $parts = $text -split ","
    $trimmedParts = @()
    foreach ($part in $parts) {
        $trimmedParts += $part.Trim()
    }

    if ($trimmedParts -notcontains "System") {
        continue
    }

    $reversedParts = @()
    for ($i = $trimmedParts.Count - 1; $i -ge 0; $i--) {
        $reversedParts += $trimmedParts[$i]
    }

    $cell.Value = [string]::Join(", ", $reversedParts)
}
